# Apply the "Add data for 2021-12-01" update to the carjacking-by-neighborhood
# by-month report: the report's as-of date moves from Nov 22 to Nov 23, and a
# handful of (mostly prior-November, i.e. late-reported) monthly counts are
# bumped up / newly populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / as-of date bookkeeping -----------------------------------
$ws.Name = "Through 2021-11-23"
$ws.Range("B1").Value = "November 2021 (through November 23)"

# --- Cell value updates -------------------------------------------------
# North Lawndale (row 2)
$ws.Range("M2").Value = 15
$ws.Range("X2").Value = 4

# Garfield Park (row 3)
$ws.Range("B3").Value = 7
$ws.Range("M3").Value = 14

# Austin (row 4)
$ws.Range("M4").Value = 10

# West Town (row 6)
$ws.Range("B6").Value = 11

# South Shore (row 7)
$ws.Range("M7").Value = 9
$ws.Range("BE7").Value = 2

# Grand Crossing (row 13) - new cell
$ws.Range("BE13").Value = 1

# Calumet Heights (row 15)
$ws.Range("BP15").Value = 2

# Chatham (row 20)
$ws.Range("M20").Value = 5
$ws.Range("X20").Value = 1

# Lake View (row 23)
$ws.Range("M23").Value = 3

# West Pullman (row 24)
$ws.Range("X24").Value = 2
$ws.Range("AI24").Value = 2

# Albany Park (row 30)
$ws.Range("M30").Value = 4

# Ashburn (row 37) - new cell
$ws.Range("AI37").Value = 1

# United Center (row 45)
$ws.Range("B45").Value = 4

# Logan Square (row 49)
$ws.Range("AT49").Value = 2

# Little Village (row 59)
$ws.Range("B59").Value = 4

# Bridgeport (row 64)
$ws.Range("M64").Value = 2

# Clearing (row 67) - new cell
$ws.Range("BE67").Value = 1

# Douglas (row 68)
$ws.Range("M68").Value = 4

# Mckinley Park (row 81) - new cell
$ws.Range("X81").Value = 1

# Portage Park (row 89) - new cell + update
$ws.Range("B89").Value = 1
$ws.Range("AI89").Value = 3
